# Update "latest output (run 136)" optimisation results.
# - Schedule sheet: recompute Cost ($) and Unit Cost ($/ML) for the two scheduled
#   pumping windows (rows 2 and 3).
# - Detailed sheet: refresh the 48-hour half-hourly price/forecast table with the
#   latest run's data (values shift forward by one period, the oldest historical
#   half-hour drops off, and Price/Pump_Status are refreshed with new forecast
#   numbers), shrinking the table from 48 rows of data (+header) to 47 rows of
#   data (+header).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Schedule sheet: update Cost ($) / Unit Cost ($/ML) for the two rows.
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 265.6614577500002
$wsSchedule.Range("F2").Value = 5.856734077380958
$wsSchedule.Range("E3").Value = 493.49900925
$wsSchedule.Range("F3").Value = 32.63882336309524

# ---------------------------------------------------------------------------
# Detailed sheet: rewrite rows 2-48 with the updated run's data, then drop the
# now-unused trailing row 49.
# ---------------------------------------------------------------------------
$wsDetailed = $wb.Worksheets.Item("Detailed")

# Each entry: RowNumber, DateTime (Excel serial), Price, Type, Pump_Status
# (column D "Date" stays 46043 for every data row, unchanged by this update)
$detailedRows = @(
    @(2, 46043.02083333334, 68.5821, "historical", "OFF"),
    @(3, 46043.04166666666, 65.96821, "historical", "OFF"),
    @(4, 46043.0625, 65.06777, "forecast", "OFF"),
    @(5, 46043.08333333334, 64.37947, "forecast", "OFF"),
    @(6, 46043.10416666666, 65, "forecast", "OFF"),
    @(7, 46043.125, 65.93982, "forecast", "OFF"),
    @(8, 46043.14583333334, 73.20005, "forecast", "OFF"),
    @(9, 46043.16666666666, 77.94, "forecast", "ON"),
    @(10, 46043.1875, 78, "forecast", "ON"),
    @(11, 46043.20833333334, 78.69043000000001, "forecast", "ON"),
    @(12, 46043.22916666666, 85.87945000000001, "forecast", "ON"),
    @(13, 46043.25, 90.55671, "forecast", "ON"),
    @(14, 46043.27083333334, 73.19, "forecast", "ON"),
    @(15, 46043.29166666666, 35.88, "forecast", "ON"),
    @(16, 46043.3125, 0.66809, "forecast", "ON"),
    @(17, 46043.33333333334, -1.08188, "forecast", "ON"),
    @(18, 46043.35416666666, -5.97577, "forecast", "ON"),
    @(19, 46043.375, -6.76009, "forecast", "ON"),
    @(20, 46043.39583333334, -7.37402, "forecast", "ON"),
    @(21, 46043.41666666666, -8.86666, "forecast", "ON"),
    @(22, 46043.4375, -9.5, "forecast", "ON"),
    @(23, 46043.45833333334, -14, "forecast", "ON"),
    @(24, 46043.47916666666, -14.96781, "forecast", "ON"),
    @(25, 46043.5, -15.94209, "forecast", "ON"),
    @(26, 46043.52083333334, -22.48859, "forecast", "ON"),
    @(27, 46043.54166666666, -22.8607, "forecast", "ON"),
    @(28, 46043.5625, -22.54286, "forecast", "ON"),
    @(29, 46043.58333333334, -23.5, "forecast", "ON"),
    @(30, 46043.60416666666, -25.33723, "forecast", "ON"),
    @(31, 46043.625, -24.13719, "forecast", "ON"),
    @(32, 46043.64583333334, -22.9965, "forecast", "ON"),
    @(33, 46043.66666666666, -6.8, "forecast", "OFF"),
    @(34, 46043.6875, -5.51, "forecast", "OFF"),
    @(35, 46043.70833333334, 36.06, "forecast", "OFF"),
    @(36, 46043.72916666666, 48.11341, "forecast", "OFF"),
    @(37, 46043.75, 55.37617, "forecast", "OFF"),
    @(38, 46043.77083333334, 61.05212, "forecast", "OFF"),
    @(39, 46043.79166666666, 73.19, "forecast", "OFF"),
    @(40, 46043.8125, 79.95, "forecast", "OFF"),
    @(41, 46043.83333333334, 73.37, "forecast", "ON"),
    @(42, 46043.85416666666, 65, "forecast", "ON"),
    @(43, 46043.875, 57.76674, "forecast", "ON"),
    @(44, 46043.89583333334, 59.86367, "forecast", "ON"),
    @(45, 46043.91666666666, 57.31, "forecast", "ON"),
    @(46, 46043.9375, 63.17447, "forecast", "ON"),
    @(47, 46043.95833333334, 64.8901, "forecast", "ON"),
    @(48, 46043.97916666666, 64.77785, "forecast", "ON")
)

foreach ($row in $detailedRows) {
    $r = $row[0]
    $wsDetailed.Cells.Item($r, 1).Value = $row[1]
    $wsDetailed.Cells.Item($r, 2).Value = $row[2]
    $wsDetailed.Cells.Item($r, 3).Value = $row[3]
    $wsDetailed.Cells.Item($r, 4).Value = 46043
    $wsDetailed.Cells.Item($r, 5).Value = $row[4]
}

# The table now only spans rows 2-48 (47 data rows); remove the stale last row.
$wsDetailed.Rows.Item(49).Delete()
